$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4114408
$ws.Range("J17").Value = 4114408
$ws.Range("L17").Value = 12343224
$ws.Range("N17").Value = -12343560
$ws.Range("H99").Value = 790.8125
$ws.Range("J99").Value = 733.3333
$ws.Range("L99").Value = 2199.9999
$ws.Range("N99").Value = -5195.9999
$ws.Range("H112").Value = 15873961
$ws.Range("J112").Value = 17858156
$ws.Range("L112").Value = 53574468
$ws.Range("N112").Value = -53576684
$ws.Range("H129").Value = 1090.0233
$ws.Range("I129").Value = 458
$ws.Range("J129").Value = 1173.1842
$ws.Range("K129").Value = 1374
$ws.Range("L129").Value = 3519.5526
$ws.Range("M129").Value = 3626
$ws.Range("N129").Value = -13519.5526
$ws.Range("H132").Value = 4341
$ws.Range("I132").Value = 3677.8572
$ws.Range("J132").Value = 5501.5
$ws.Range("K132").Value = 11033.5716
$ws.Range("L132").Value = 16504.5
$ws.Range("M132").Value = -8503.571599999999
$ws.Range("N132").Value = -21564.5
$ws.Range("H138").Value = 1758.6268
$ws.Range("I138").Value = 1102.1482
$ws.Range("J138").Value = 2201.75
$ws.Range("K138").Value = 3306.4446
$ws.Range("L138").Value = 6605.25
$ws.Range("M138").Value = 1833.5554
$ws.Range("N138").Value = -16885.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5710.476
$ws.Range("I45").Value = 6682.706
$ws.Range("J45").Value = 1578.5
$ws.Range("K45").Value = 6682.706
$ws.Range("L45").Value = 1578.5
$ws.Range("M45").Value = -6305.706
$ws.Range("N45").Value = -2332.5
$ws.Range("H61").Value = 273235.9
$ws.Range("I61").Value = 2209.7778
$ws.Range("K61").Value = 2209.7778
$ws.Range("M61").Value = -1997.7778
$ws.Range("H63").Value = 125003624
$ws.Range("I63").Value = 125003624
$ws.Range("K63").Value = 125003624
$ws.Range("M63").Value = -125002938
$ws.Range("H66").Value = 125003624
$ws.Range("I66").Value = 125003624
$ws.Range("K66").Value = 625018120
$ws.Range("M66").Value = -625014688
$ws.Range("H74").Value = 1641.625
$ws.Range("I74").Value = 1417.7368
$ws.Range("K74").Value = 1417.7368
$ws.Range("M74").Value = -543.7367999999999
$ws.Range("H77").Value = 1641.625
$ws.Range("I77").Value = 1417.7368
$ws.Range("K77").Value = 7088.683999999999
$ws.Range("M77").Value = -2720.683999999999
$ws.Range("H88").Value = 2573.889
$ws.Range("I88").Value = 2354.2
$ws.Range("J88").Value = 2848.5
$ws.Range("K88").Value = 2354.2
$ws.Range("L88").Value = 2848.5
$ws.Range("M88").Value = -1948.2
$ws.Range("N88").Value = -3660.5
$ws.Range("H91").Value = 2573.889
$ws.Range("I91").Value = 2354.2
$ws.Range("J91").Value = 2848.5
$ws.Range("K91").Value = 2354.2
$ws.Range("L91").Value = 2848.5
$ws.Range("M91").Value = -950.1999999999998
$ws.Range("N91").Value = -5656.5
$ws.Range("H122").Value = 1767.7667
$ws.Range("I122").Value = 1781.4828
$ws.Range("K122").Value = 5344.4484
$ws.Range("M122").Value = -2894.4484
$ws.Range("H136").Value = 273235.9
$ws.Range("I136").Value = 2209.7778
$ws.Range("K136").Value = 6629.3334
$ws.Range("M136").Value = -4079.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11655.208
$ws.Range("I20").Value = 1091.75
$ws.Range("J20").Value = 32782.125
$ws.Range("K20").Value = 1091.75
$ws.Range("L20").Value = 32782.125
$ws.Range("M20").Value = -844.75
$ws.Range("N20").Value = -33276.125
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H86").Value = 1642
$ws.Range("I86").Value = 1383.25
$ws.Range("J86").Value = 2332
$ws.Range("K86").Value = 1383.25
$ws.Range("L86").Value = 2332
$ws.Range("M86").Value = -260.25
$ws.Range("N86").Value = -4578
$ws.Range("H89").Value = 1642
$ws.Range("I89").Value = 1383.25
$ws.Range("J89").Value = 2332
$ws.Range("K89").Value = 6916.25
$ws.Range("L89").Value = 11660
$ws.Range("M89").Value = -1300.25
$ws.Range("N89").Value = -22892
$ws.Range("H99").Value = 1184.6522
$ws.Range("I99").Value = 782.7143
$ws.Range("K99").Value = 782.7143
$ws.Range("M99").Value = 715.2857
$ws.Range("H105").Value = 22605.3
$ws.Range("I105").Value = 27514.375
$ws.Range("K105").Value = 27514.375
$ws.Range("M105").Value = -25767.375
$ws.Range("H134").Value = 2547.8928
$ws.Range("I134").Value = 2218.6428
$ws.Range("K134").Value = 6655.928400000001
$ws.Range("M134").Value = -4120.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1675.3784
$ws.Range("I58").Value = 1489.2106
$ws.Range("J58").Value = 1871.8889
$ws.Range("K58").Value = 1489.2106
$ws.Range("L58").Value = 1871.8889
$ws.Range("M58").Value = -1286.2106
$ws.Range("N58").Value = -2277.8889
$ws.Range("H94").Value = 3678.081
$ws.Range("I94").Value = 3527.8572
$ws.Range("J94").Value = 3769.5217
$ws.Range("K94").Value = 3527.8572
$ws.Range("L94").Value = 3769.5217
$ws.Range("M94").Value = -3076.8572
$ws.Range("N94").Value = -4671.521699999999
$ws.Range("H122").Value = 1854514.5
$ws.Range("I122").Value = 2778667.5
$ws.Range("K122").Value = 8336002.5
$ws.Range("M122").Value = -8333552.5
$ws.Range("H134").Value = 3447.138
$ws.Range("I134").Value = 3424.7036
$ws.Range("K134").Value = 10274.1108
$ws.Range("M134").Value = -7739.110799999999
$ws.Range("H136").Value = 1675.3784
$ws.Range("I136").Value = 1489.2106
$ws.Range("J136").Value = 1871.8889
$ws.Range("K136").Value = 4467.6318
$ws.Range("L136").Value = 5615.6667
$ws.Range("M136").Value = -1917.6318
$ws.Range("N136").Value = -10715.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5463.485
$ws.Range("I70").Value = 5453.4287
$ws.Range("J70").Value = 5519.8
$ws.Range("K70").Value = 5453.4287
$ws.Range("L70").Value = 5519.8
$ws.Range("M70").Value = -5183.4287
$ws.Range("N70").Value = -6059.8
$ws.Range("H73").Value = 5463.485
$ws.Range("I73").Value = 5453.4287
$ws.Range("J73").Value = 5519.8
$ws.Range("K73").Value = 5453.4287
$ws.Range("L73").Value = 5519.8
$ws.Range("M73").Value = -4517.4287
$ws.Range("N73").Value = -7391.8
$ws.Range("H80").Value = 5488.2256
$ws.Range("I80").Value = 6678.591
$ws.Range("J80").Value = 2578.4443
$ws.Range("K80").Value = 6678.591
$ws.Range("L80").Value = 2578.4443
$ws.Range("M80").Value = -5680.591
$ws.Range("N80").Value = -4574.4443
$ws.Range("H83").Value = 5488.2256
$ws.Range("I83").Value = 6678.591
$ws.Range("J83").Value = 2578.4443
$ws.Range("K83").Value = 33392.955
$ws.Range("L83").Value = 12892.2215
$ws.Range("M83").Value = -28400.955
$ws.Range("N83").Value = -22876.2215
$ws.Range("H102").Value = 362049.44
$ws.Range("I102").Value = 678929.6
$ws.Range("K102").Value = 678929.6
$ws.Range("M102").Value = -677307.6
$ws.Range("H113").Value = 1718.0454
$ws.Range("I113").Value = 1163.3636
$ws.Range("J113").Value = 2272.7273
$ws.Range("K113").Value = 1163.3636
$ws.Range("L113").Value = 2272.7273
$ws.Range("M113").Value = 1006.6364
$ws.Range("N113").Value = -6612.7273
$ws.Range("H126").Value = 4109.2744
$ws.Range("I126").Value = 5811.8696
$ws.Range("J126").Value = 2710.7144
$ws.Range("K126").Value = 17435.6088
$ws.Range("L126").Value = 8132.1432
$ws.Range("M126").Value = -14965.6088
$ws.Range("N126").Value = -13072.1432
$ws.Range("H132").Value = 3049.303
$ws.Range("I132").Value = 3194.7896
$ws.Range("K132").Value = 9584.3688
$ws.Range("M132").Value = -7054.3688

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2181.9092
$ws.Range("I40").Value = 2289.6843
$ws.Range("J40").Value = 1499.3334
$ws.Range("K40").Value = 2289.6843
$ws.Range("L40").Value = 1499.3334
$ws.Range("M40").Value = -2153.6843
$ws.Range("N40").Value = -1771.3334
$ws.Range("H122").Value = 503877.2
$ws.Range("I122").Value = 4030.2222
$ws.Range("K122").Value = 12090.6666
$ws.Range("M122").Value = -9640.6666
$ws.Range("H136").Value = 5007.4287
$ws.Range("I136").Value = 2306.5
$ws.Range("J136").Value = 13650.4
$ws.Range("K136").Value = 6919.5
$ws.Range("L136").Value = 40951.2
$ws.Range("M136").Value = -4369.5
$ws.Range("N136").Value = -46051.2
$ws.Range("H140").Value = 45824
$ws.Range("J140").Value = 45824
$ws.Range("L140").Value = 45824
$ws.Range("N140").Value = -56184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10018.333
$ws.Range("J11").Value = 15002.5
$ws.Range("L11").Value = 15002.5
$ws.Range("N11").Value = -15286.5
$ws.Range("H14").Value = 70005
$ws.Range("J14").Value = 70005
$ws.Range("L14").Value = 70005
$ws.Range("N14").Value = -70341
$ws.Range("H46").Value = 39000
$ws.Range("J46").Value = 39000
$ws.Range("L46").Value = 39000
$ws.Range("N46").Value = -39462
$ws.Range("H132").Value = 1311.0769
$ws.Range("I132").Value = 931.3946999999999
$ws.Range("J132").Value = 2341.6428
$ws.Range("K132").Value = 2794.1841
$ws.Range("L132").Value = 7024.928400000001
$ws.Range("M132").Value = -264.1840999999999
$ws.Range("N132").Value = -12084.9284
$ws.Range("H134").Value = 39000
$ws.Range("J134").Value = 39000
$ws.Range("L134").Value = 117000
$ws.Range("N134").Value = -122070
